# Generate Report for Handback
# The a1611c47-fe89-4e11-a2a2-904b596ac550 file has now been handed back
# (it was previously "Ready for handoff"). Update the Overview sheet and
# the per-locale (zh-cn / de-de) status sheets to reflect the handback:
#  - Status flips from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Handback DateTime is refreshed
#  - Error Detail (stale handback-version warning) is cleared

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-09-05 06:55:36"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-09-05 06:55:45"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839
